$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting/style used by the existing header cells (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Column I is a constant "1" for every data row (2-29)
$ws.Range("I2:I29").Value = 1

# Column J duplicates the values already present in column H (2-29)
$ws.Range("H2:H29").Copy()
$ws.Range("J2:J29").PasteSpecial(-4163)

$excel.CutCopyMode = $false
